$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 39 (shifts existing rows 39..100 down to 40..101)
$ws.Rows(39).Insert()

# Populate the new row 39 with this week's record (same Fruta/Coco template as
# surrounding rows, new date + volume/price figures)
$ws.Range("A39").Value2 = 6
$ws.Range("B39").Value2 = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C39").Value2 = "Metropolitana"
$ws.Range("D39").Value2 = 45259
$ws.Range("E39").Value2 = 13
$ws.Range("F39").Value2 = "Fruta"
$ws.Range("G39").Value2 = 100108
$ws.Range("H39").Value2 = "Tropicales y subtropicales"
$ws.Range("I39").Value2 = 100108007
$ws.Range("J39").Value2 = "Coco"
$ws.Range("K39").Value2 = "Sin especificar"
$ws.Range("L39").Value2 = "Primera"
$ws.Range("M39").Value2 = 100
$ws.Range("N39").Value2 = 28000
$ws.Range("O39").Value2 = 28000
$ws.Range("P39").Value2 = 28000
$ws.Range("Q39").Value2 = "$/malla 20 unidades"
$ws.Range("R39").Value2 = "Perú"
$ws.Range("S39").Value2 = 1400
$ws.Range("T39").Value2 = 20
